$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy number-format/style from existing cells of the same column so the
# shared cellXfs entries (s="1", s="4", s="2") are reused instead of new
# style entries being minted.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("C1").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)

$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# New row of data: #7, 2025-07-02, "Instalação Geladeira", 250, "Araujo"
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = 45840
$ws.Range("C10").Value = "Instalação Geladeira"
$ws.Range("D10").Value = 250
$ws.Range("E10").Value = "Araujo"

$ws.Range("E10").Select()
